# Update "想去人数" (F column) counts on three sheets of the workbook.
# Sheet "展览" (rId1 / sheet1)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 387
$ws1.Range("F9").Value = 1233
$ws1.Range("F11").Value = 282
$ws1.Range("F14").Value = 6841
$ws1.Range("F18").Value = 7750
$ws1.Range("F21").Value = 4489
$ws1.Range("F23").Value = 2251
$ws1.Range("F24").Value = 955
$ws1.Range("F26").Value = 238
$ws1.Range("F30").Value = 263
$ws1.Range("F31").Value = 225
$ws1.Range("F33").Value = 1935
$ws1.Range("F41").Value = 2040

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 248

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 248
$ws4.Range("F11").Value = 1233
$ws4.Range("F13").Value = 282
$ws4.Range("F16").Value = 6841
$ws4.Range("F20").Value = 7750
$ws4.Range("F23").Value = 4489
$ws4.Range("F25").Value = 2251
$ws4.Range("F26").Value = 955
$ws4.Range("F28").Value = 238
$ws4.Range("F34").Value = 263
$ws4.Range("F36").Value = 1935
$ws4.Range("F45").Value = 2040
